# Add season record columns (Wins, Losses, Ties) to the player stats sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: copy formatting from an existing header cell (AC1) so the new
# header cells pick up the same bold/border/alignment style, then overwrite
# the copied text with the real header labels.
$ws.Range("AC1").Copy($ws.Range("AD1"))
$ws.Range("AD1").Value = "Wins"

$ws.Range("AC1").Copy($ws.Range("AE1"))
$ws.Range("AE1").Value = "Losses"

$ws.Range("AC1").Copy($ws.Range("AF1"))
$ws.Range("AF1").Value = "Ties"

# Data rows 2-50: every row gets the same season record.
for ($row = 2; $row -le 50; $row++) {
    $ws.Cells.Item($row, 30).Value = 84
    $ws.Cells.Item($row, 31).Value = 78
    $ws.Cells.Item($row, 32).Value = 0
}
